# Auto update Excel log
# Appends newly captured sensor log entries to the "Proximity" and "Camera"
# worksheets of the SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$RowNum,
        [string]$Date,
        [string]$Timestamp,
        [string]$Hour,
        [string]$Location,
        [string]$Value,
        [string]$Status
    )

    # Column A holds dates formatted as plain "yyyy-mm-dd" text in this log.
    # Force the cell to Text first so the engine doesn't reinterpret the
    # string as a real date serial value, then restore the default
    # "Normal" style so no stray formatting is left behind on the cell.
    $ws.Range("A$RowNum").NumberFormat = "@"
    $ws.Range("A$RowNum").Value = $Date
    $ws.Range("A$RowNum").Style = "Normal"

    $ws.Range("B$RowNum").Value = $Timestamp
    $ws.Range("C$RowNum").Value = $Hour
    $ws.Range("D$RowNum").Value = $Location
    $ws.Range("E$RowNum").Value = $Value
    $ws.Range("F$RowNum").Value = $Status
}

# --- Proximity sheet: append new ENTER/EXIT events as rows 55-57 ---
$wsProximity = $wb.Worksheets.Item("Proximity")

Add-LogRow $wsProximity 55 "2026-02-01" "14:47:02" "14:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $wsProximity 56 "2026-02-01" "14:47:23" "14:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $wsProximity 57 "2026-02-01" "14:47:26" "14:00" "Living Room Main Door" "EXIT"  "User EXITED Living Room Main Door"

# --- Camera sheet: append new Image Received/Captured events as rows 38-39 ---
$wsCamera = $wb.Worksheets.Item("Camera")

Add-LogRow $wsCamera 38 "2026-02-01" "14:47:02" "14:00" "Living Room Main Door" "Image Received" "Active"
Add-LogRow $wsCamera 39 "2026-02-01" "14:47:26" "14:00" "Living Room Main Door" "Image Captured" "Active"
